# "productos importar, solo campos necesarios"
# The import only needs a single generic "ID" column up front instead of
# the three separate ID columns (ID tipo, ID categoría, ID Kit) that used
# to be scattered through the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column for the "ID" field.
$ws.Range("A1").EntireColumn.Insert()
$ws.Range("A1").Value = "ID"

# Drop the columns the import doesn't need. Work right-to-left so the
# earlier column letters stay valid while we go.
$ws.Range("Q1").EntireColumn.Delete()
$ws.Range("E1").EntireColumn.Delete()
$ws.Range("C1").EntireColumn.Delete()

# Give the new "ID" header the same look as the other short "code" header
# (Referencia, now back in column C) already uses.
$ws.Range("C1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reset the view back to the top-left corner instead of the stale R2
# selection left over from the deleted "ID s Instituciones" lookup.
[void]$ws.Range("A1").Select()
